$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the new column header first so "x_i" becomes shared-string index 0.
$ws.Range("B1").Value = "x_i"

# Column A: keep the General number format but centre the whole extended
# range (A1:A25) -- the sheet grows downward to give the accompanying
# chart more rows to plot against.
$ws.Range("A1:A25").NumberFormat = "General"
$ws.Range("A1:A25").HorizontalAlignment = -4108   # xlCenter

# B1 (the header) shares the same look as column A.
$ws.Range("B1").NumberFormat = "General"
$ws.Range("B1").HorizontalAlignment = -4108   # xlCenter

# Row 10 gets a text label, "9 iteracija" (last iteration marker).
$ws.Range("A10").Value = "9 iteracija"

# New column B: a shifted-down copy of column A's original values, each
# with a fixed 6-decimal numeric format, centred.
$ws.Range("B2").Value = 5
$ws.Range("B3").Value = 3.3783783783783785
$ws.Range("B4").Value = 2.320008885077228
$ws.Range("B5").Value = 1.6487812280309717
$ws.Range("B6").Value = 1.2528031537556461
$ws.Range("B7").Value = 1.0604121514975322
$ws.Range("B8").Value = 1.0047989143538689
$ws.Range("B9").Value = 1.0000341618854527
$ws.Range("B10").Value = 1.000000001750412

$ws.Range("B2:B10").NumberFormat = "0.000000"
$ws.Range("B2:B10").HorizontalAlignment = -4108   # xlCenter

# Widen the columns to fit the new data, per the authored layout.
$ws.Columns.Item(1).ColumnWidth = 8.140625
$ws.Columns.Item(2).ColumnWidth = 14.41

$ws.Range("E18").Select() | Out-Null
